$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.2639766931533813, 0.9125000238418579),
    @(0.2397506684064865, 0.9262499809265137),
    @(0.3285242021083832, 0.8862500190734863),
    @(0.317320704460144,  0.887499988079071),
    @(0.3031081259250641, 0.8987500071525574),
    @(0.2830060124397278, 0.9162499904632568)
)

$startRow = 8
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$wb.Save()
